# Resuelto: exceso de creditos
#
# The two signature blocks ("{coordinador}" and "{jefe}") are each followed
# by a paragraph that only contains the underline
# "____________________________________". Those paragraphs only carried a
# Complex Script font (w:cs="Arial"); the Western/ASCII font was left to
# fall back to the theme's default (Cambria), so the underline printed in
# the wrong font. This sets the Western font (w:ascii/w:hAnsi) to Arial on
# both underline paragraphs while keeping the existing Complex Script Arial
# font intact.

$d = $word.ActiveDocument

function Find-ExactTagParagraph($doc, $tag) {
    # Locate the (single) paragraph whose whole text is exactly $tag, e.g.
    # "{coordinador}" or "{jefe}" -- ignoring the trailing paragraph mark.
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $para = $doc.Paragraphs($i)
        $text = $para.Range.Text.TrimEnd([char]13, [char]7)
        if ($text -eq $tag) {
            return $para
        }
    }
    return $null
}

# --- Underline below "{coordinador}" ---
# Apply to the whole paragraph (run + paragraph mark), which preserves the
# paragraph mark's existing w:cs="Arial" while adding w:ascii/w:hAnsi.
$coordTag  = Find-ExactTagParagraph $d "{coordinador}"
$coordLine = $coordTag.Next()
$coordLine.Range.Font.Name = "Arial"

# --- Underline below "{jefe}" ---
# Apply only to the run text (not the paragraph mark), using Find to select
# just the underline characters, then set both the Western font and
# re-assert the Complex Script font so the existing w:cs="Arial" survives.
$jefeTag  = Find-ExactTagParagraph $d "{jefe}"
$jefeLine = $jefeTag.Next()
$jefeRange = $jefeLine.Range
$jefeRange.Find.ClearFormatting()
$jefeRange.Find.Execute("____________________________________", $false, $false, $false, $false, $false, $true, 1, $false, "", 0, $false)
$jefeRange.Font.Name = "Arial"
$jefeRange.Font.NameBi = "Arial"
